$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1 = "is" (no special fill)
$ws.Range("B1").Value = "is"

# C1 = "good" with green fill
$ws.Range("C1").Value = "good"
$ws.Range("C1").Interior.Color = 32768

# Row 2: A2:C2 empty numeric cells with yellow fill
$ws.Range("A2:C2").Interior.Color = 65535
